$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slides 3 and 6 swap their entire visual content.
#    Slide 3 currently holds the "ANALYSIS OF STOCK OF ELECTRIC VEHICLES
#    BY COUNTRY" picture + two text boxes; Slide 6 currently holds a
#    single full-bleed "Non-Tesla EV Models" dashboard picture.
#    After the edit, slide 3 shows the dashboard picture (resized) and
#    slide 6 shows the EV-by-country picture + the two text boxes.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s6 = $p.Slides.Item(6)

# Remember how many shapes originally belong to each slide before we
# start appending copies.
$s3OrigCount = $s3.Shapes.Count
$s6OrigCount = $s6.Shapes.Count

# Copy slide 6's original picture onto slide 3 first.
$s6.Shapes.Item(1).Copy()
$movedToS3 = $p.Slides.Item(3).Shapes.Paste()

# Copy slide 3's original shapes (picture + 2 textboxes, in order) onto
# slide 6.
for ($i = 1; $i -le $s3OrigCount; $i++) {
    $sh = $s3.Shapes.Item($i)
    $sh.Copy()
    [void]$s6.Shapes.Paste()
}

# Remove the old shapes that were copied from (leave only the migrated
# ones behind).
for ($i = 1; $i -le $s3OrigCount; $i++) {
    $s3.Shapes.Item(1).Delete()
}
for ($i = 1; $i -le $s6OrigCount; $i++) {
    $s6.Shapes.Item(1).Delete()
}

# Slide 3 now only contains the migrated dashboard picture; nudge it to
# its final, slightly-retouched position/size.
$pic3 = $s3.Shapes.Item(1)
$pic3.Left = 1630016 / 12700
$pic3.Top = 172278 / 12700
$pic3.Width = 10243932 / 12700
$pic3.Height = 5685183 / 12700

# ---------------------------------------------------------------------
# 2) Slide 4: retitle the title placeholder, reposition/resize it, and
#    move the supporting screenshot.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$title4 = $s4.Shapes.Item(1)
$title4.Left = 516835 / 12700
$title4.Top = 1934817 / 12700
$title4.Width = 3578087 / 12700
$title4.Height = 1137567 / 12700
$title4.TextFrame.TextRange.Text = "Assessing the relationship between the size of the stores, number of employees and revenue? "

$pic4 = $s4.Shapes.Item(3)
$pic4.Left = 4232945 / 12700
$pic4.Top = 583096 / 12700
